# Apply the "interactor instructions" notes sheet + highlight a few budget rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("begroting")

# Highlight (apply the existing "Good" cell style, index 16) on B9, B20, B21
# to match B7/B10/B11/... which already use it.
$ws1.Range("B9").Style = "Good"
$ws1.Range("B20").Style = "Good"
$ws1.Range("B21").Style = "Good"

# Update selection on begroting sheet: activeCell B7, selection B7:B21.
$ws1.Range("B7:B21").Select()

# Add a new worksheet "Sheet1" right after "begroting" with interactor
# instructions / notes. It becomes the active sheet.
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Sheet1"

$ws2.Range("A1").Value = "ctrl + x"
$ws2.Range("B1").Value = "interactor instructions"
$ws2.Range("A2").Value = "shift + c "
$ws2.Range("B2").Value = "polygong offset"
$ws2.Range("A3").Value = "use buffers for nodes and edges"

$ws2.Range("A4").Select()

$wb.Save()
